# AHC-20-001_DP_Switch_Barrier_Block_Interface_bom.xlsx — REV D update
# Adds a new BOM line (Item 5: LCBSB-4-01A-RT board support / snap-lock
# standoff, qty 4, refs J7 J8 J9 J10) as row 6 of the BOM table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New BOM row (Item Number, Quantity, Part Reference, MFR, MFRPN,
#     Package, Value, Voltage, Description) ---------------------------
$ws.Range("A6").NumberFormat = "0"
$ws.Range("B6").NumberFormat = "0"
$ws.Range("C6:I6").NumberFormat = "@"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 4
$ws.Range("C6").Value = "J7 J8 J9 J10"
$ws.Range("D6").Value = "Essentra Components"
$ws.Range("E6").Value = "LCBSB-4-01A-RT"
$ws.Range("F6").Value = "700mil Square"
$ws.Range("G6").Value = "LCBSB-4-01A-RT"
$ws.Range("H6").Value = ""
$ws.Range("I6").Value = "BRD SPT SNAP LOCK ADHESIVE 1/4`""

# --- Column widths re-fit to the new (wider) content ------------------
$ws.Columns.Item(4).ColumnWidth = 19.6666666666667
$ws.Columns.Item(6).ColumnWidth = 12.6666666666667
$ws.Columns.Item(9).ColumnWidth = 31.6666666666667

# --- Selection left where the author's cursor ended up ----------------
$ws.Range("F9").Select() | Out-Null
